# Auto-generated script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.821.20"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.105.52"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.58"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.31"
$ws.Range("E7").Value = "  +3.05%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +2.35%  "

$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.80"
$ws.Range("E12").Value = "  +6.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.417.30"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.07"
$ws.Range("E14").Value = "  -1.62%  "

$ws.Range("E15").Value = "  +3.02%  "

$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.107.65"
$ws.Range("E17").Value = "  +0.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.840.11"
$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.73"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  +1.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.40"
$ws.Range("E22").Value = "  +1.47%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("E24").Value = "  -3.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.81"
$ws.Range("E26").Value = "  +4.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.84"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("E29").Value = "  +3.71%  "

$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  +12.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.60"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.23"
$ws.Range("E34").Value = "  +13.41%  "

$ws.Range("E35").Value = "  -1.29%  "

$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.02"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("E42").Value = "  +3.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.524.87"
$ws.Range("E43").Value = "  -0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +8.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.81"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.84"
$ws.Range("E46").Value = "  +2.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.10"

$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.17"
$ws.Range("E49").Value = "  +1.37%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.305.14"
$ws.Range("E51").Value = "  +1.40%  "
